$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 name/link change (TrustWalletToken replacing BitcoinSV)
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"

# Price (D) and Volume(1h) (E) updates
$ws.Range("D2").Value = "42.399.91"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "2.288.57"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'300.83"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").Value = "'94.93"
$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("D7").Value = "'0.507"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "  -2.12%  "

$ws.Range("D10").Value = "'34.25"
$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("D11").Value = "'19.05"
$ws.Range("E11").Value = "  +2.80%  "

$ws.Range("D12").Value = "'0.0780"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D14").Value = "'6.70"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "2.646.86"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").Value = "2.298.01"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "42.353.80"
$ws.Range("E18").Value = "  -0.89%  "

$ws.Range("D19").Value = "'12.14"
$ws.Range("E19").Value = "  -6.81%  "

$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("D21").Value = "'5.95"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "'67.43"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = "  +5.99%  "

$ws.Range("D24").Value = "'234.91"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -2.45%  "

$ws.Range("D27").Value = "'24.19"
$ws.Range("E27").Value = "  -3.72%  "

$ws.Range("D28").Value = "'2.36"
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").Value = "'164.92"
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("D30").Value = "'9.02"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").Value = "'31.62"
$ws.Range("E31").Value = "  -4.32%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'4.96"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").Value = "'17.50"
$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("D35").Value = "'0.0693"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D37").Value = "'4.32"
$ws.Range("E37").Value = "  -9.66%  "

$ws.Range("D38").Value = "'0.0996"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D39").Value = "'1.73"
$ws.Range("E39").Value = "  -1.15%  "

$ws.Range("E40").Value = "  -1.18%  "

$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("D42").Value = "'19.88"
$ws.Range("E42").Value = "  +8.84%  "

$ws.Range("D43").Value = "1.949.42"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("E44").Value = "  +2.75%  "

$ws.Range("D45").Value = "'0.0277"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").Value = "'2.10"
$ws.Range("E46").Value = "  +2.75%  "

$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -2.59%  "

$ws.Range("D48").Value = "'2.84"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "2.515.66"
$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("D50").Value = "'52.82"
$ws.Range("E50").Value = "  -1.84%  "

$ws.Range("D51").Value = "'1.13"
$ws.Range("E51").Value = "  +0.09%  "

